# Update New York GDP figures (rows 12-30) with revised data, and append
# a new observation row (31) for 2020-01-01.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised GDP values for existing observations (2001-2019) ---
$ws.Range("B12").Value = 943455.21499999997
$ws.Range("B13").Value = 956748.49899999995
$ws.Range("B14").Value = 978374.46299999999
$ws.Range("B15").Value = 1030208.245
$ws.Range("B16").Value = 1087191.9180000001
$ws.Range("B17").Value = 1152952.0149999999
$ws.Range("B18").Value = 1200662.0789999999
$ws.Range("B19").Value = 1200903.861
$ws.Range("B20").Value = 1228112.524
$ws.Range("B21").Value = 1288303.129
$ws.Range("B22").Value = 1312974.304
$ws.Range("B23").Value = 1400779.267
$ws.Range("B24").Value = 1445252.2320000001
$ws.Range("B25").Value = 1507782.7860000001
$ws.Range("B26").Value = 1570332.7450000001
$ws.Range("B27").Value = 1638128.9410000001
$ws.Range("B28").Value = 1690244.554
$ws.Range("B29").Value = 1790858.0789999999
$ws.Range("B30").Value = 1872165.5049999999

# --- Append new observation row 31 (2020-01-01) ---
# Copy formatting (date format on A, 0.000 number format on B) from the
# row above so the new cells pick up the same styles as the rest of the
# table, then fill in the values.
$ws.Range("A30:B30").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A31").Value = 43831
$ws.Range("B31").Value = 1809323.3970000001

# --- Update the visible selection to match the author's saved state ---
$ws.Range("A1:B1048576").Select() | Out-Null
